$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '22.478.83'
$ws.Range("E2").Value = '  -0.05%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.571.10'
$ws.Range("E3").Value = '  -0.20%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("E5").Value = '  +0.03%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '286.34'
$ws.Range("E6").Value = '  -2.01%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3650'
$ws.Range("E7").Value = '  -1.98%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '48.13'
$ws.Range("E8").Value = '  -3.54%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3331'

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.127'
$ws.Range("E10").Value = '  -2.18%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07430'
$ws.Range("E11").Value = '  -1.67%  '

$ws.Range("E12").Value = '  +0.15%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.87'
$ws.Range("E13").Value = '  -2.02%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.971'
$ws.Range("E14").Value = '  -1.43%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.916'
$ws.Range("E15").Value = '  -0.83%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.567.43'
$ws.Range("E16").Value = '  -0.80%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001107'
$ws.Range("E17").Value = '  -1.84%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '87.89'
$ws.Range("E18").Value = '  -3.70%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06737'
$ws.Range("E19").Value = '  -0.37%  '

$ws.Range("E20").Value = '  +0.05%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.383'
$ws.Range("E21").Value = '  +1.03%  '

$ws.Range("E22").Value = '  +0.43%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.04'
$ws.Range("E23").Value = '  -1.07%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '22.463.29'
$ws.Range("E24").Value = '  -0.14%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.381'
$ws.Range("E25").Value = '  +0.22%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.632'
$ws.Range("E26").Value = '  -1.50%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '151.54'
$ws.Range("E27").Value = '  +1.68%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.54'
$ws.Range("E28").Value = '  -2.69%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.018'
$ws.Range("E29").Value = '  -0.98%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '124.52'
$ws.Range("E30").Value = '  -0.92%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.747.01'
$ws.Range("E31").Value = '  -0.35%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.032'
$ws.Range("E32").Value = '  -4.26%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.139'
$ws.Range("E33").Value = '  -1.55%  '

$ws.Range("E34").Value = '  -1.10%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.773'
$ws.Range("E35").Value = '  -0.98%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.08262'
$ws.Range("E36").Value = '  -1.64%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02425'
$ws.Range("E37").Value = '  -2.73%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2245'
$ws.Range("E38").Value = '  -2.74%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06453'
$ws.Range("E39").Value = '  -1.35%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.415'

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.294'
$ws.Range("E41").Value = '  -3.07%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '11.31'
$ws.Range("E42").Value = '  -0.60%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.6279'
$ws.Range("E43").Value = '  +0.41%  '

$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.89'
$ws.Range("E44").Value = '  -1.22%  '

$ws.Range("B45").Value = 'Decentraland'
$ws.Range("C45").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6091'
$ws.Range("E45").Value = '  +4.25%  '

$ws.Range("B46").Value = 'PancakeSwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.748'
$ws.Range("E46").Value = '  -1.70%  '

$ws.Range("B47").Value = 'NEARProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.047'
$ws.Range("E47").Value = '  -1.96%  '

$ws.Range("B48").Value = 'Quant'
$ws.Range("C48").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '124.25'
$ws.Range("E48").Value = '  -4.77%  '

$ws.Range("B49").Value = 'EOS'
$ws.Range("C49").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.215'
$ws.Range("E49").Value = '  -0.65%  '

$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07212'
$ws.Range("E50").Value = '  -1.68%  '

$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '76.35'
$ws.Range("E51").Value = '  -0.66%  '

